$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r`a") -eq "Food (Seeds/Chocolate/Peppers)") {
        $p.Range.Delete()
        break
    }
}
